# Apply the edit described by the diff:
#  - parameters!B3  (fleet_size) : 4 -> 7
#  - parameters!B12 (no_stops)   : 6 -> 7
#  - comp_quantity_inst1: remove the "T3 -> T4, 4, 1" row entirely (row 6),
#    which shifts the remaining rows up; then the four rows that used to read
#    T1->T2 / T2->T3 / T3->T4 / T4->T5 (quantity, 0) are re-pointed so their
#    origin becomes "F1" (and the first one's destination becomes T1 as well)
#  - comp_quantity_inst1!C5 (was C5=28 in the old layout, the T2->T3 "required"
#    row) becomes 1
#  - selections are moved to match the saved view (parameters!D4,
#    comp_quantity_inst1!C9) with comp_quantity_inst1 left as the active sheet/tab

$wb = $excel.ActiveWorkbook

# --- Sheet: parameters ---------------------------------------------------
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("B3").Value = 7
$wsParams.Range("B12").Value = 7

# --- Sheet: comp_quantity_inst1 ------------------------------------------
$wsComp = $wb.Worksheets.Item("comp_quantity_inst1")

# Row 5 (T2 -> T3, required) quantity drops from 28 to 1
$wsComp.Range("C5").Value = 1

# Remove the "T3 -> T4" required row outright; rows below shift up one place
$wsComp.Rows(6).Delete() | Out-Null

# The (now) D=0 block's origins all become "F1"; the first row's destination
# also changes from T2 to T1
$wsComp.Range("A6").Value = "F1"
$wsComp.Range("B6").Value = "T1"
$wsComp.Range("A7").Value = "F1"
$wsComp.Range("A8").Value = "F1"
$wsComp.Range("A9").Value = "F1"

# --- Selections / active view -------------------------------------------
# Touch parameters' selection first (any sheet visited before the final one
# loses "active" status once we select on comp_quantity_inst1 last).
$wsParams.Range("D4").Select() | Out-Null

$wsTrip = $wb.Worksheets.Item("trip_duration")
$wsTrip.Range("D21").Select() | Out-Null

# comp_quantity_inst1 stays the active/visible tab, matching the saved file.
$wsComp.Range("C9").Select() | Out-Null
